$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.797.25"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.493.18"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'593.22"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'172.27"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  +3.95%  "
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "'0.430"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").Value = "4.098.47"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "'29.43"
$ws.Range("E14").Value = "  +4.91%  "
$ws.Range("D15").Value = "66.826.42"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "3.445.72"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "'6.26"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "'14.19"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").Value = "'393.13"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'7.90"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'73.28"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'0.533"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'10.17"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'0.993"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "'6.12"
$ws.Range("E29").Value = "  -2.91%  "
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "'23.57"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "'7.33"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'1.60"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "'162.64"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").Value = "'6.84"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").Value = "'4.62"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "2.834.94"
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'27.13"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0736"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "'26.04"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "'42.65"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").Value = "'2.52"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "'0.0301"
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("D47").Value = "'337.41"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").Value = "'34.43"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "'6.41"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'0.839"
$ws.Range("E51").Value = "  -3.19%  "

# Reset style on text-forced numeric-looking cells so they keep
# the sheet default style (no explicit quote-prefix style index).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
